$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 100
$ws.Range("B3").Value = "EFECTIVO"
$ws.Range("A4").Value = 20000
$ws.Range("B4").Value = "DEBITO"
